$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "Yes"
$ws.Range("D3").Value = "No"
$ws.Range("D4").Value = "Yes"
